# Updated symbol list on Fri Dec 23 14:28:24 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.16"
$ws.Range("D3").Value = "'21.98"
$ws.Range("D4").Value = "'5.405"
$ws.Range("D7").Value = "'6.359"
$ws.Range("D8").Value = "'0.8175"
$ws.Range("D9").Value = "'1.013"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("D11").Value = "'0.07465"
$ws.Range("D12").Value = "'0.03452"
$ws.Range("D13").Value = "'0.03040"
$ws.Range("D14").Value = "'4.216"
$ws.Range("D15").Value = "'0.09393"
$ws.Range("D16").Value = "'0.001609"
$ws.Range("D17").Value = "'0.04831"
$ws.Range("D19").Value = "'0.006019"
$ws.Range("D20").Value = "'0.004098"
$ws.Range("D21").Value = "'0.001001"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D23").Value = "'3.697"
$ws.Range("D24").Value = "'2.220"
$ws.Range("D26").Value = "'0.1297"
$ws.Range("D40").Value = "'0.03863"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002414"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003013"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.006158"
$ws.Range("D47").Value = "'0.8808"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").Value = "'0.1423"
